$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 137 (ALC)
$ws.Range("H137").Value = 1637.3778
$ws.Range("I137").Value = 1315.9166
$ws.Range("J137").Value = 1754.2727
$ws.Range("K137").Value = 3947.7498
$ws.Range("L137").Value = 5262.8181
$ws.Range("M137").Value = -1397.7498
$ws.Range("N137").Value = -10362.8181

# Row 138 (ALC)
$ws.Range("H138").Value = 4499.8335
$ws.Range("I138").Value = 2158.4
$ws.Range("J138").Value = 5877.147
$ws.Range("K138").Value = 6475.200000000001
$ws.Range("L138").Value = 17631.441
$ws.Range("M138").Value = -1335.200000000001
$ws.Range("N138").Value = -27911.441

$ws = $wb.Worksheets.Item("ARM")
# Row 3 (ARM)
$ws.Range("H3").Value = 1470
$ws.Range("I3").Value = 1470
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1470
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -1355
$ws.Range("N3").ClearContents()

# Row 74 (ARM)
$ws.Range("H74").Value = 1939.4773
$ws.Range("I74").Value = 1416.8695
$ws.Range("J74").Value = 2511.8572
$ws.Range("K74").Value = 1416.8695
$ws.Range("L74").Value = 2511.8572
$ws.Range("M74").Value = -542.8695
$ws.Range("N74").Value = -4259.8572

# Row 77 (ARM)
$ws.Range("H77").Value = 1939.4773
$ws.Range("I77").Value = 1416.8695
$ws.Range("J77").Value = 2511.8572
$ws.Range("K77").Value = 7084.3475
$ws.Range("L77").Value = 12559.286
$ws.Range("M77").Value = -2716.3475
$ws.Range("N77").Value = -21295.286

# Row 122 (ARM)
$ws.Range("H122").Value = 1849.579
$ws.Range("I122").Value = 1849.579
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5548.737
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3098.737
$ws.Range("N122").ClearContents()

# Row 123 (ARM)
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

# Row 124 (ARM)
$ws.Range("H124").Value = 34164.5
$ws.Range("J124").Value = 34164.5
$ws.Range("L124").Value = 34164.5
$ws.Range("N124").Value = -43984.5

# Row 125 (ARM)
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

# Row 127 (ARM)
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

# Row 132 (ARM)
$ws.Range("H132").Value = 3158.423
$ws.Range("I132").Value = 3406.325
$ws.Range("J132").Value = 2332.0833
$ws.Range("K132").Value = 10218.975
$ws.Range("L132").Value = 6996.249899999999
$ws.Range("M132").Value = -7688.974999999999
$ws.Range("N132").Value = -12056.2499

$ws = $wb.Worksheets.Item("BSM")
# Row 57 (BSM)
$ws.Range("H57").Value = 38000
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()

# Row 135 (BSM)
$ws.Range("H135").Value = 68000
$ws.Range("J135").Value = 68000
$ws.Range("L135").Value = 68000
$ws.Range("N135").Value = -78140

# Row 136 (BSM)
$ws.Range("H136").Value = 38000
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 4 (CRP)
$ws.Range("H4").Value = 1500
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()

# Row 31 (CRP)
$ws.Range("H31").Value = 835.13
$ws.Range("J31").Value = 777.2655999999999
$ws.Range("L31").Value = 777.2655999999999
$ws.Range("N31").Value = -1367.2656

# Row 34 (CRP)
$ws.Range("H34").Value = 835.13
$ws.Range("J34").Value = 777.2655999999999
$ws.Range("L34").Value = 777.2655999999999
$ws.Range("N34").Value = -1181.2656

# Row 54 (CRP)
$ws.Range("H54").Value = 15092
$ws.Range("J54").Value = 15092
$ws.Range("L54").Value = 15092
$ws.Range("N54").Value = -16408

# Row 108 (CRP)
$ws.Range("H108").Value = 27781
$ws.Range("J108").Value = 27781
$ws.Range("L108").Value = 27781
$ws.Range("N108").Value = -35461

# Row 122 (CRP)
$ws.Range("H122").Value = 2304.6897
$ws.Range("I122").Value = 2110.182
$ws.Range("J122").Value = 2916
$ws.Range("K122").Value = 6330.545999999999
$ws.Range("L122").Value = 8748
$ws.Range("M122").Value = -3880.545999999999
$ws.Range("N122").Value = -13648

# Row 134 (CRP)
$ws.Range("H134").Value = 2251.0557
$ws.Range("I134").Value = 1527.6666
$ws.Range("J134").Value = 2974.4443
$ws.Range("K134").Value = 4582.9998
$ws.Range("L134").Value = 8923.332900000001
$ws.Range("M134").Value = -2047.9998
$ws.Range("N134").Value = -13993.3329

$ws = $wb.Worksheets.Item("CUL")
# Row 131 (CUL)
$ws.Range("H131").Value = 1504.6111
$ws.Range("J131").Value = 1537.8049
$ws.Range("L131").Value = 4613.4147
$ws.Range("N131").Value = -14693.4147

$ws = $wb.Worksheets.Item("GSM")
# Row 5 (GSM)
$ws.Range("H5").Value = 1263575
$ws.Range("I5").Value = 2500500
$ws.Range("J5").Value = 26650
$ws.Range("K5").Value = 2500500
$ws.Range("L5").Value = 26650
$ws.Range("M5").Value = -2500388
$ws.Range("N5").Value = -26874

# Row 24 (GSM)
$ws.Range("H24").Value = 2692.3076
$ws.Range("J24").Value = 2692.3076
$ws.Range("L24").Value = 2692.3076
$ws.Range("N24").Value = -3038.3076

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (LTW)
$ws.Range("H7").Value = 1611.3334
$ws.Range("I7").Value = 1365.875
$ws.Range("J7").Value = 3575
$ws.Range("K7").Value = 1365.875
$ws.Range("L7").Value = 3575
$ws.Range("M7").Value = -1253.875
$ws.Range("N7").Value = -3799

# Row 126 (LTW)
$ws.Range("H126").Value = 1611.3334
$ws.Range("I126").Value = 1365.875
$ws.Range("J126").Value = 3575
$ws.Range("K126").Value = 4097.625
$ws.Range("L126").Value = 10725
$ws.Range("M126").Value = -1627.625
$ws.Range("N126").Value = -15665

$ws = $wb.Worksheets.Item("WVR")
# Row 122 (WVR)
$ws.Range("H122").Value = 1494
$ws.Range("I122").Value = 1494
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4482
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2032
$ws.Range("N122").ClearContents()

# Row 132 (WVR)
$ws.Range("H132").Value = 2874.75
$ws.Range("I132").Value = 3285.3928
$ws.Range("J132").Value = 2299.85
$ws.Range("K132").Value = 9856.178400000001
$ws.Range("L132").Value = 6899.549999999999
$ws.Range("M132").Value = -7326.178400000001
$ws.Range("N132").Value = -11959.55

# Row 136 (WVR)
$ws.Range("H136").Value = 1024.0714
$ws.Range("I136").Value = 723.1177
$ws.Range("J136").Value = 1489.1818
$ws.Range("K136").Value = 2169.3531
$ws.Range("L136").Value = 4467.5454
$ws.Range("M136").Value = 380.6468999999997
$ws.Range("N136").Value = -9567.545399999999
